# SMARTNODES.docx (Chinese Traditional) - apply translation-commit edit:
#   1. Reset the "smarthosting" bookmark's id from 1 to 0.
#   2. Remove the "CryptoBridge" / "HitBTC" HYPERLINK field codes (and the
#      now-unneeded "such as" / ", " wording around them) from the
#      "100,000 Smart:" paragraph, leaving "... obtained from exchanges."

$d = $word.ActiveDocument

# --- Change 1: bookmark "smarthosting" id 1 -> 0 --------------------------
# The Word object model does not expose bookmark ids directly, but deleting
# a bookmark and re-adding it over the same range causes Word to renumber
# it starting at 0 (the lowest free id), matching the diff's id="0".
$bm = $d.Bookmarks("smarthosting")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("smarthosting", $bmRange) | Out-Null

# --- Change 2: drop the CryptoBridge / HitBTC exchange links --------------
# Fields(1) = HYPERLINK "https://crypto-bridge.org/" -> "CryptoBridge"
# Fields(2) = HYPERLINK "https://hitbtc.com/SMART-to-BTC" -> "HitBTC"
$f1 = $d.Fields.Item(1)
$f2 = $d.Fields.Item(2)

# The fldChar "begin" marker for field 1 sits exactly one character before
# its Code range starts; the " such as" text (plus its trailing non-
# breaking space) ends exactly there, 9 characters earlier.
$fldBeginPos = $f1.Code.Start - 1
$suchAsStart = $fldBeginPos - 9

# Remove the ", " separator that sits between field 1's result and field 2's
# code (i.e. between "CryptoBridge" and the HitBTC field).
$d.Range($f1.Result.End, $f2.Code.Start).Delete()

# Delete the fields themselves (highest offset first so earlier positions
# stay valid). Field.Delete() removes the begin/separate/end fldChar runs,
# the instrText run and the displayed result text in one go.
$f2.Delete()
$f1.Delete()

# Finally trim the now-orphaned " such as" (+ trailing nbsp) that used to
# introduce the removed links, leaving "... from exchanges." directly.
$d.Range($suchAsStart, $fldBeginPos).Delete()
